$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates: force Text format to preserve exact formatting ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.624.95"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.685.55"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.12"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5342"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2680"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06439"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.70"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07812"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.694.83"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.512"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5615"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8441"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.03"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.667.12"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.795"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "196.00"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.376"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.003"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.40"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1281"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.483"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.35"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.435"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06159"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.612"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.471"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.702"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.799"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5755"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01646"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.042"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.071.93"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8638"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.45"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.837.13"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.31"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.189"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05215"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.081"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +3.41%  "
$ws.Range("E3").Value = "  +4.46%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("E6").Value = "  +3.81%  "
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("E8").Value = "  +5.47%  "
$ws.Range("E9").Value = "  +4.82%  "
$ws.Range("E10").Value = "  +9.23%  "
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("E13").Value = "  +4.91%  "
$ws.Range("E14").Value = "  +4.76%  "
$ws.Range("E15").Value = "  +8.38%  "
$ws.Range("E16").Value = "  +3.56%  "
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("E20").Value = "  +7.62%  "
$ws.Range("E21").Value = "  +5.30%  "
$ws.Range("E22").Value = "  +6.50%  "
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("E25").Value = "  +7.11%  "
$ws.Range("E26").Value = "  +3.38%  "
$ws.Range("E27").Value = "  +6.80%  "
$ws.Range("E28").Value = "  +5.43%  "
$ws.Range("E29").Value = "  +5.54%  "
$ws.Range("E30").Value = "  +3.65%  "
$ws.Range("E31").Value = "  +9.42%  "
$ws.Range("E32").Value = "  +5.79%  "
$ws.Range("E33").Value = "  +8.35%  "
$ws.Range("E34").Value = "  +6.31%  "
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("E36").Value = "  +1.41%  "
$ws.Range("E37").Value = "  +1.74%  "
$ws.Range("E38").Value = "  +4.77%  "
$ws.Range("E39").Value = "  +8.41%  "
$ws.Range("E40").Value = "  +6.06%  "
$ws.Range("E41").Value = "  +4.43%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("E44").Value = "  +3.18%  "
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("E46").Value = "  +7.33%  "
$ws.Range("E47").Value = "  +4.36%  "
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("E50").Value = "  +6.67%  "

# --- Row 51: Mantle -> RenderToken (name, link, price, volume all change) ---
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.470"
$ws.Range("E51").Value = "  +9.09%  "
